$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 2659
$ws.Range("I3").Value = 2790
$ws.Range("I4").Value = 678
$ws.Range("G5").Value = 785
$ws.Range("I5").Value = 247
$ws.Range("I6").Value = 3173
$ws.Range("G7").Value = 24650
$ws.Range("I7").Value = 9547

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I4").Value = 37
$ws.Range("I6").Value = 67
$ws.Range("I7").Value = 316
$ws.Range("I8").Value = 609
$ws.Range("I14").Value = 47
$ws.Range("I15").Value = 123
$ws.Range("I17").Value = 10
$ws.Range("I19").Value = 265
$ws.Range("I20").Value = 240
$ws.Range("I23").Value = 86
$ws.Range("I29").Value = 638
$ws.Range("I30").Value = 29
$ws.Range("I33").Value = 450
$ws.Range("I34").Value = 39
$ws.Range("I35").Value = 14
$ws.Range("I36").Value = 127
$ws.Range("I37").Value = 311
$ws.Range("I42").Value = 330
$ws.Range("I43").Value = 86
$ws.Range("I48").Value = 104
$ws.Range("I49").Value = 64
$ws.Range("I51").Value = 84
$ws.Range("I52").Value = 196
$ws.Range("I53").Value = 105
$ws.Range("I55").Value = 105
$ws.Range("I57").Value = 33
$ws.Range("I59").Value = 19
$ws.Range("I61").Value = 12
$ws.Range("G63").Value = 196
$ws.Range("I63").Value = 35
$ws.Range("I64").Value = 86
$ws.Range("I65").Value = 214
$ws.Range("I67").Value = 364
$ws.Range("I76").Value = 150
$ws.Range("I77").Value = 51
$ws.Range("I78").Value = 129
$ws.Range("I79").Value = 245
$ws.Range("I83").Value = 188
$ws.Range("I84").Value = 75
$ws.Range("I85").Value = 442
$ws.Range("I89").Value = 103
$ws.Range("I90").Value = 108
$ws.Range("I91").Value = 113
$ws.Range("I94").Value = 84
$ws.Range("I96").Value = 116
$ws.Range("G101").Value = 24650
$ws.Range("I101").Value = 9547

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I3").Value = 182
$ws.Range("I5").Value = 14
$ws.Range("I6").Value = 116
$ws.Range("I7").Value = 442

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I2").Value = 49
$ws.Range("I7").Value = 196

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 193
$ws.Range("I3").Value = 168
$ws.Range("I4").Value = 38
$ws.Range("I6").Value = 191
$ws.Range("I7").Value = 609

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I2").Value = 18
$ws.Range("I7").Value = 105

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I3").Value = 96
$ws.Range("I6").Value = 80
$ws.Range("I7").Value = 316

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("I6").Value = 36
$ws.Range("I7").Value = 103

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I2").Value = 31
$ws.Range("I7").Value = 116

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("I3").Value = 12
$ws.Range("I7").Value = 47

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("I2").Value = 8
$ws.Range("I7").Value = 29

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 101
$ws.Range("I7").Value = 311

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I3").Value = 127
$ws.Range("I7").Value = 364

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("I2").Value = 30
$ws.Range("I7").Value = 75

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 66
$ws.Range("I3").Value = 57
$ws.Range("I6").Value = 69
$ws.Range("I7").Value = 214

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I2").Value = 67
$ws.Range("I7").Value = 188

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 103
$ws.Range("I3").Value = 165
$ws.Range("I7").Value = 450

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("I2").Value = 17
$ws.Range("I6").Value = 35
$ws.Range("I7").Value = 64

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 198
$ws.Range("I3").Value = 223
$ws.Range("I5").Value = 23
$ws.Range("I6").Value = 173
$ws.Range("I7").Value = 638

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I3").Value = 71
$ws.Range("I7").Value = 265

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I3").Value = 21
$ws.Range("I7").Value = 104

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I2").Value = 32
$ws.Range("I3").Value = 38
$ws.Range("I7").Value = 150

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I2").Value = 29
$ws.Range("I6").Value = 12
$ws.Range("I7").Value = 67

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I2").Value = 86
$ws.Range("I7").Value = 330

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I4").Value = 19
$ws.Range("I7").Value = 129

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I2").Value = 36
$ws.Range("I3").Value = 28
$ws.Range("I6").Value = 35
$ws.Range("I7").Value = 105

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I2").Value = 23
$ws.Range("I7").Value = 86

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I3").Value = 37
$ws.Range("I7").Value = 113

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 69
$ws.Range("I3").Value = 78
$ws.Range("I7").Value = 245

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("I6").Value = 31
$ws.Range("I7").Value = 86

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I2").Value = 63
$ws.Range("I6").Value = 88
$ws.Range("I7").Value = 240

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("I2").Value = 3
$ws.Range("I7").Value = 10

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I2").Value = 39
$ws.Range("I3").Value = 39
$ws.Range("I7").Value = 127

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("I2").Value = 14
$ws.Range("I7").Value = 39

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I6").Value = 48
$ws.Range("I7").Value = 84

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I2").Value = 38
$ws.Range("I3").Value = 28
$ws.Range("I7").Value = 123

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("I6").Value = 5
$ws.Range("I7").Value = 14

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("I4").Value = 2

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range("I7").Value = 19

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I4").Value = 12
$ws.Range("I7").Value = 108

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I2").Value = 15
$ws.Range("I3").Value = 23
$ws.Range("I7").Value = 84

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("I2").Value = 12
$ws.Range("I4").Value = 6
$ws.Range("I7").Value = 33

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("I2").Value = 16
$ws.Range("I6").Value = 49
$ws.Range("I7").Value = 86

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("I4").Value = 1
$ws.Range("I7").Value = 51

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("I4").Value = 3
$ws.Range("I7").Value = 37

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("I2").Value = 4

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range("I7").Value = 12
